# This workbook is a monthly budget tracker. Each worksheet (except "result")
# lists 12 month-blocks; each block has an "Allocation" row, an "Expend" row
# and a "Remaining" row (in that order) whose values live in column B.
#
# The commit "change Analyze file to month/month" clears out the budget
# figures for every month except the first one on most sheets, and on the
# two sheets that track a single lump-sum commitment ("دانشگاه"/University
# and "پس انداز"/Savings) it shifts the Expend/Remaining split from the
# first month into an adjusted split, zeroing all of the middle months and
# finally turning the Remaining of the last month negative.

$wb = $excel.ActiveWorkbook

# --- Sheets 1 & 2: Allocation, Expend AND Remaining all reset to 0
#     for every month after the first one. ---
$sheetsABC = @(1, 2)
$allocRowsABC = @{
    1 = @(6, 14, 22, 30, 38, 46, 54, 62, 70, 78, 86, 94)
    2 = @(5, 13, 21, 29, 37, 45, 53, 61, 69, 77, 85, 93)
}
foreach ($sheetIdx in $sheetsABC) {
    $ws = $wb.Worksheets.Item($sheetIdx)
    $rows = $allocRowsABC[$sheetIdx]
    for ($m = 1; $m -lt $rows.Count; $m++) {
        $allocRow = $rows[$m]
        $ws.Range("B$allocRow").Value = 0
        $ws.Range("B$($allocRow + 1)").Value = 0
        $ws.Range("B$($allocRow + 2)").Value = 0
    }
}

# --- Sheets 3, 4, 5, 6 & 9: Expend is already 0 on these sheets, so only
#     Allocation and Remaining need to be reset to 0 for months 2-12. ---
$sheetsAR = @(3, 4, 5, 6, 9)
$allocRowsAR = @(4, 12, 20, 28, 36, 44, 52, 60, 68, 76, 84, 92)
foreach ($sheetIdx in $sheetsAR) {
    $ws = $wb.Worksheets.Item($sheetIdx)
    for ($m = 1; $m -lt $allocRowsAR.Count; $m++) {
        $allocRow = $allocRowsAR[$m]
        $ws.Range("B$allocRow").Value = 0
        $ws.Range("B$($allocRow + 2)").Value = 0
    }
}

# --- Sheets 7 (دانشگاه) & 8 (پس انداز): single lump-sum commitments.
#     Month 1: Expend -> 0, Remaining -> the original allocation amount.
#     Months 2-11: Allocation and Expend -> 0 (Remaining already 0).
#     Month 12: Allocation -> 0, Remaining -> minus the original allocation
#               amount (Expend is left untouched). ---
$lumpSum = @{
    7 = @{ Rows = @(4, 12, 20, 28, 36, 44, 52, 60, 68, 76, 84, 93); Amount = 1400000 }
    8 = @{ Rows = @(4, 12, 20, 28, 36, 44, 52, 60, 68, 76, 84, 94); Amount = 600000 }
}
foreach ($sheetIdx in @(7, 8)) {
    $ws = $wb.Worksheets.Item($sheetIdx)
    $rows = $lumpSum[$sheetIdx].Rows
    $amount = $lumpSum[$sheetIdx].Amount

    # Month 1 (index 0)
    $firstAlloc = $rows[0]
    $ws.Range("B$($firstAlloc + 1)").Value = 0
    $ws.Range("B$($firstAlloc + 2)").Value = $amount

    # Months 2-11 (indices 1..count-2)
    for ($m = 1; $m -lt ($rows.Count - 1); $m++) {
        $allocRow = $rows[$m]
        $ws.Range("B$allocRow").Value = 0
        $ws.Range("B$($allocRow + 1)").Value = 0
    }

    # Month 12 (last index)
    $lastAlloc = $rows[$rows.Count - 1]
    $ws.Range("B$lastAlloc").Value = 0
    $ws.Range("B$($lastAlloc + 2)").Value = (0 - $amount)
}
